# Update the header date line.
$d = $word.ActiveDocument
$d.Content.Find.Execute("2025-11-16 Sunday", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "2025-11-17 Monday", 2)

# Update each answer cell in the single table via explicit (row, column)
# addressing. This avoids ambiguity since some old values repeat
# (e.g. "25÷8=3, 1" appears twice but maps to two different new values).
$t = $d.Tables.Item(1)

$updates = @(
    @{ Row = 1;  Col = 1; New = "62÷8=7, 6" },
    @{ Row = 1;  Col = 2; New = "64÷3=21, 1" },
    @{ Row = 1;  Col = 3; New = "22÷3=7, 1" },
    @{ Row = 1;  Col = 4; New = "38÷5=7, 3" },
    @{ Row = 1;  Col = 5; New = "33÷9=3, 6" },

    @{ Row = 5;  Col = 1; New = "24÷9=2, 6" },
    @{ Row = 5;  Col = 2; New = "94÷6=15, 4" },
    @{ Row = 5;  Col = 3; New = "48÷4=12, 0" },
    @{ Row = 5;  Col = 4; New = "56÷9=6, 2" },
    @{ Row = 5;  Col = 5; New = "42÷5=8, 2" },

    @{ Row = 9;  Col = 1; New = "36÷2=18, 0" },
    @{ Row = 9;  Col = 2; New = "90÷9=10, 0" },
    @{ Row = 9;  Col = 3; New = "18÷2=9, 0" },
    @{ Row = 9;  Col = 4; New = "47÷9=5, 2" },
    @{ Row = 9;  Col = 5; New = "80÷9=8, 8" },

    @{ Row = 13; Col = 1; New = "83÷5=16, 3" },
    @{ Row = 13; Col = 2; New = "40÷8=5, 0" },
    @{ Row = 13; Col = 3; New = "27÷6=4, 3" },
    @{ Row = 13; Col = 4; New = "12÷9=1, 3" },
    @{ Row = 13; Col = 5; New = "35÷2=17, 1" },

    @{ Row = 17; Col = 1; New = "69÷4=17, 1" },
    @{ Row = 17; Col = 2; New = "33÷7=4, 5" },
    @{ Row = 17; Col = 3; New = "75÷6=12, 3" },
    @{ Row = 17; Col = 4; New = "61÷4=15, 1" },
    @{ Row = 17; Col = 5; New = "82÷4=20, 2" }
)

foreach ($u in $updates) {
    $cell = $t.Cell($u.Row, $u.Col)
    $cell.Range.Text = $u.New
}
